# actualizacion de base de datos 320 y quitar el tiempo de espera entre numeros 3seg
#
# Replace the phone-number list in column A with a new 10-number rotation
# (without the "57" country-code prefix), shrink the list from 351 down to
# 320 rows, and clear out a stray "-" that was typed into B10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new set of 10 numbers, in the order they should first appear
# (this also controls the order they land in the shared-strings table).
$numbers = @(
    "3215996243",
    "3204530013",
    "3157520347",
    "3134688382",
    "3105694409",
    "3104023154",
    "3174466432",
    "3183978799",
    "3183247990",
    "3016406749"
)

$totalRows = 320

for ($r = 1; $r -le $totalRows; $r++) {
    $idx = ($r - 1) % $numbers.Length
    $ws.Cells.Item($r, 1).Value = $numbers[$idx]
}

# Stray manual edit that slipped into B10 during the cleanup.
$ws.Cells.Item(10, 2).Value = "-"

# Drop the now-unused tail rows (321-351) so the sheet shrinks to 320 rows.
$ws.Range("A321:A351").EntireRow.Delete()

# Restore the view: scrolled near the bottom of the (now shorter) list,
# with C7 selected.
$excel.ActiveWindow.ScrollRow = 312
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C7").Select() | Out-Null
